$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5 from
# 45207 (2023-10-08) to 45208 (2023-10-09), matching the source diff.
$ws.Range("C2:C5").Value = 45208
